$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet is a daily price log, newest entry on row 2. A new day's
# entry (16-01-2026) is being prepended, which pushes every existing
# row down by one position (the previous final row, row 219, becomes
# row 220). Column F holds hyperlinked PDF links for the most recent
# ~158 rows; since the runtime does not shift hyperlink anchors when
# rows are inserted, we drop them first and recreate them afterwards
# against their new row numbers.
# ------------------------------------------------------------------

$ws.Hyperlinks.Delete()

$ws.Rows("2:2").Insert()

# The freshly inserted row inherits bold formatting from the header
# above it; reset it to the plain "data row" look used everywhere
# else in the sheet before writing the new values into it.
for ($c = 1; $c -le 6; $c++) {
    $cell = $ws.Cells.Item(2, $c)
    $cell.Style = "Normal"
    $cell.Font.Bold = $false
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}
$ws.Range("D2").NumberFormat = "0.000"

$ws.Range("A2").Value = "16-01-2026"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 343.25
$ws.Range("E2").Value = "15.01.2026"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-15-january-2026.pdf"

# Recreate the hyperlinks for column F (rows 2 through 159 now carry
# a PDF link, since every row shifted down by one and the previous
# last-linked row, 158, is now row 159).
for ($r = 2; $r -le 159; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value2
    if ($url) {
        $ws.Hyperlinks.Add($cell, $url)
        $cell.Style = "Normal"
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
    }
}
